# The assignment Q&A answering "What are the further improvements that can
# be made to make it efficient?" gets two more bullet points appended to its
# list (same numbered/bulleted list as the existing "We can use python
# programming language ..." item): one about key eviction, one about
# persisting session data to a file.

$d = $word.ActiveDocument

# Locate the paragraph that ends the existing bullet list ("We can use
# python programming language ...") by matching its distinctive text
# rather than hard-coding an index, so the script is resilient to minor
# structural differences.
$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*We can use python programming language*") {
        $anchorIndex = $i
        break
    }
}

# Insert the first new bullet right after the anchor paragraph. Using
# InsertParagraphAfter on the paragraph's Range clones the paragraph's
# list/indent/justification formatting (numPr numId=4, ilvl=0, ind
# left=1440 hanging=360, jc=both) exactly like pressing Enter at the end of
# that bullet in Word would. Re-fetch paragraphs by index after each
# mutation rather than caching object references, since the collection
# reflows on every insert.
$anchorPara = $d.Paragraphs.Item($anchorIndex)
$anchorPara.Range.InsertParagraphAfter()
$bullet1Index = $anchorIndex + 1
$bullet1 = $d.Paragraphs.Item($bullet1Index)
$bullet1.Range.Text = "In the provided implementation there is no way to do the eviction of the keys which are expired. As it would be a problem if we have lot of expire key."

# Insert the second new bullet right after the first one, same list.
$bullet1 = $d.Paragraphs.Item($bullet1Index)
$bullet1.Range.InsertParagraphAfter()
$bullet2Index = $bullet1Index + 1
$bullet2 = $d.Paragraphs.Item($bullet2Index)
$bullet2.Range.Text = "The given data for a particular session can be store in some file and a function can be constructed to import the data when the somehow the connection is lost."
